$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row 3: CCAC动漫游戏嘉年华, row 10: AEO纯白礼赞动漫嘉年华
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 103
$wsExhibit.Range("F10").Value = 406

# Sheet "全部类型" (all types, mirrors the same data) - same two rows/cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 103
$wsAll.Range("F10").Value = 406

$wb.Save()
